$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.233.20"
$ws.Range("E2").Value = "  +3.21%  "
$ws.Range("D3").Value = "3.061.91"
$ws.Range("E3").Value = "  +2.71%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.447"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.46%  "
$ws.Range("E9").Value = "  +2.97%  "
$ws.Range("E10").Value = "  +5.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.370"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.48%  "
$ws.Range("D12").Value = "3.542.28"
$ws.Range("E12").Value = "  +1.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.130"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.67%  "
$ws.Range("E15").Value = "  +14.91%  "
$ws.Range("D16").Value = "58.200.28"
$ws.Range("E16").Value = "  +3.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +9.38%  "
$ws.Range("D18").Value = "3.071.66"
$ws.Range("E18").Value = "  +3.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "338.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.27%  "
$ws.Range("E22").Value = "  +1.20%  "
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.503"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.170"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.66%  "
$ws.Range("D27").Value = "0.0₃0962"
$ws.Range("E27").Value = "  +8.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.95"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +10.37%  "
$ws.Range("E31").Value = "  +5.85%  "
$ws.Range("E32").Value = "  +4.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.61%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.80"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.20%  "
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "157.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.74%  "
$ws.Range("E37").Value = "  +2.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.49"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +10.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0696"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.45%  "
$ws.Range("D40").Value = "3.100.71"
$ws.Range("E40").Value = "  +2.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "37.82"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.04%  "
$ws.Range("E42").Value = "  +9.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.668"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.35%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").Value = "2.336.56"
$ws.Range("E45").Value = "  +4.69%  "
$ws.Range("E46").Value = "  +4.87%  "
$ws.Range("E47").Value = "  +3.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0242"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.95%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.87"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.89"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.47%  "
